$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "SE151252"
$ws.Range("C2").Value = "BAP GROUP"
$ws.Range("D2").Value = "Kỹ thuật phần mềm"
$ws.Range("E2").Value = 6.0
$ws.Range("F2").Value = "Nhiệt huyết trong công việc."
